# Scheduled-runner style update: refresh cached market price / profit
# figures on several rows across multiple job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW). Values below are plain numeric overwrites (no formulas
# are used in this workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 780
$ws.Range("I98").Value = 764.75
$ws.Range("J98").Value = 902
$ws.Range("K98").Value = 764.75
$ws.Range("L98").Value = 902
$ws.Range("M98").Value = 733.25
$ws.Range("N98").Value = -3898

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 780
$ws.Range("I122").Value = 764.75
$ws.Range("J122").Value = 902
$ws.Range("K122").Value = 2294.25
$ws.Range("L122").Value = 2706
$ws.Range("M122").Value = 155.75
$ws.Range("N122").Value = -7606

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2463.1924
$ws.Range("I61").Value = 2195.4119
$ws.Range("K61").Value = 2195.4119
$ws.Range("M61").Value = -1983.4119

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1847.4
$ws.Range("I74").Value = 1864.6364
$ws.Range("K74").Value = 1864.6364
$ws.Range("M74").Value = -990.6364000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1847.4
$ws.Range("I77").Value = 1864.6364
$ws.Range("K77").Value = 9323.182000000001
$ws.Range("M77").Value = -4955.182000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2897
$ws.Range("I132").Value = 1727
$ws.Range("J132").Value = 7723.25
$ws.Range("K132").Value = 5181
$ws.Range("L132").Value = 23169.75
$ws.Range("M132").Value = -2651
$ws.Range("N132").Value = -28229.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2463.1924
$ws.Range("I136").Value = 2195.4119
$ws.Range("K136").Value = 6586.2357
$ws.Range("M136").Value = -4036.2357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1746.8235
$ws.Range("I20").Value = 1710.5454
$ws.Range("J20").Value = 1813.3334
$ws.Range("K20").Value = 1710.5454
$ws.Range("L20").Value = 1813.3334
$ws.Range("M20").Value = -1463.5454
$ws.Range("N20").Value = -2307.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8335232.5
$ws.Range("I86").Value = 10754435
$ws.Range("K86").Value = 10754435
$ws.Range("M86").Value = -10753312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 8335232.5
$ws.Range("I89").Value = 10754435
$ws.Range("K89").Value = 53772175
$ws.Range("M89").Value = -53766559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1712.5161
$ws.Range("I94").Value = 1418.5454
$ws.Range("J94").Value = 2431.111
$ws.Range("K94").Value = 1418.5454
$ws.Range("L94").Value = 2431.111
$ws.Range("M94").Value = -967.5454
$ws.Range("N94").Value = -3333.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2703.318
$ws.Range("I31").Value = 1357.0465
$ws.Range("K31").Value = 1357.0465
$ws.Range("M31").Value = -1062.0465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2703.318
$ws.Range("I34").Value = 1357.0465
$ws.Range("K34").Value = 1357.0465
$ws.Range("M34").Value = -1155.0465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 972.8461
$ws.Range("I122").Value = 904.1177
$ws.Range("K122").Value = 2712.3531
$ws.Range("M122").Value = -262.3531000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 5882383
$ws.Range("I38").Value = 8333354
$ws.Range("J38").Value = 52.4
$ws.Range("K38").Value = 25000062
$ws.Range("L38").Value = 157.2
$ws.Range("M38").Value = -24999715
$ws.Range("N38").Value = -851.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 71428710
$ws.Range("I75").Value = 293
$ws.Range("K75").Value = 879
$ws.Range("M75").Value = 119

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 71428710
$ws.Range("I78").Value = 293
$ws.Range("K78").Value = 2637
$ws.Range("M78").Value = 2355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 9000
$ws.Range("N84").Value = -20232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 509.54544
$ws.Range("I86").Value = 493.66666
$ws.Range("J86").Value = 528.6
$ws.Range("K86").Value = 1480.99998
$ws.Range("L86").Value = 1585.8
$ws.Range("M86").Value = -294.9999800000001
$ws.Range("N86").Value = -3957.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2000
$ws.Range("I87").Value = 2000
$ws.Range("K87").Value = 6000
$ws.Range("M87").Value = -4752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 509.54544
$ws.Range("I89").Value = 493.66666
$ws.Range("J89").Value = 528.6
$ws.Range("K89").Value = 4442.99994
$ws.Range("L89").Value = 4757.400000000001
$ws.Range("M89").Value = 1485.00006
$ws.Range("N89").Value = -16613.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 2000
$ws.Range("I90").Value = 2000
$ws.Range("K90").Value = 18000
$ws.Range("M90").Value = -11760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1057.6786
$ws.Range("I121").Value = 656
$ws.Range("J121").Value = 1145
$ws.Range("K121").Value = 1968
$ws.Range("L121").Value = 3435
$ws.Range("M121").Value = -658
$ws.Range("N121").Value = -6055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 5849.7856
$ws.Range("I124").Value = 5000
$ws.Range("J124").Value = 5915.154
$ws.Range("K124").Value = 15000
$ws.Range("L124").Value = 17745.462
$ws.Range("M124").Value = -10090
$ws.Range("N124").Value = -27565.462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3031233.5
$ws.Range("I131").Value = 11111437
$ws.Range("J131").Value = 1157.4584
$ws.Range("K131").Value = 33334311
$ws.Range("L131").Value = 3472.3752
$ws.Range("M131").Value = -33329271
$ws.Range("N131").Value = -13552.3752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 15880819
$ws.Range("I137").Value = 2104.2856
$ws.Range("J137").Value = 23820176
$ws.Range("K137").Value = 6312.8568
$ws.Range("L137").Value = 71460528
$ws.Range("M137").Value = -1212.8568
$ws.Range("N137").Value = -71470728

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2156.75
$ws.Range("J140").Value = 2000
$ws.Range("L140").Value = 6000
$ws.Range("N140").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1721
$ws.Range("I102").Value = 1616
$ws.Range("J102").Value = 1966
$ws.Range("K102").Value = 1616
$ws.Range("L102").Value = 1966
$ws.Range("M102").Value = 6
$ws.Range("N102").Value = -5210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 31250900
$ws.Range("I93").Value = 957.1429000000001
$ws.Range("K93").Value = 957.1429000000001
$ws.Range("M93").Value = 290.8570999999999
